# Refresh market-price derived columns (H-N) across the Leve profit
# sheets, per the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1857.7273
$ws.Range("J17").Value = 1857.7273
$ws.Range("L17").Value = 5573.1819
$ws.Range("N17").Value = -5909.1819
$ws.Range("H70").Value = 10149.75
$ws.Range("I70").Value = 10839.6
$ws.Range("K70").Value = 32518.8
$ws.Range("M70").Value = -32248.8
$ws.Range("H73").Value = 10149.75
$ws.Range("I73").Value = 10839.6
$ws.Range("K73").Value = 32518.8
$ws.Range("M73").Value = -31582.8
$ws.Range("H100").Value = 1500
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = ""
$ws.Range("H141").Value = 2177.7144
$ws.Range("I141").Value = 1395.625
$ws.Range("J141").Value = 3220.5
$ws.Range("K141").Value = 4186.875
$ws.Range("L141").Value = 9661.5
$ws.Range("M141").Value = 993.125
$ws.Range("N141").Value = -20021.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3670122
$ws.Range("I32").Value = 3503628
$ws.Range("K32").Value = 3503628
$ws.Range("M32").Value = -3503341
$ws.Range("H132").Value = 2123.5
$ws.Range("I132").Value = 1999
$ws.Range("K132").Value = 5997
$ws.Range("M132").Value = -3467
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4320.1665
$ws.Range("I86").Value = 4320.1665
$ws.Range("K86").Value = 4320.1665
$ws.Range("M86").Value = -3197.1665
$ws.Range("H89").Value = 4320.1665
$ws.Range("I89").Value = 4320.1665
$ws.Range("K89").Value = 21600.8325
$ws.Range("M89").Value = -15984.8325
$ws.Range("H94").Value = 2006.7778
$ws.Range("I94").Value = 1951.75
$ws.Range("J94").Value = 2447
$ws.Range("K94").Value = 1951.75
$ws.Range("L94").Value = 2447
$ws.Range("M94").Value = -1500.75
$ws.Range("N94").Value = -3349
$ws.Range("H99").Value = 2529.875
$ws.Range("I99").Value = 2391.2856
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 2391.2856
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -893.2856000000002
$ws.Range("N99").Value = -6496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 199141.5
$ws.Range("J9").Value = 199141.5
$ws.Range("L9").Value = 199141.5
$ws.Range("N9").Value = -199477.5
$ws.Range("H22").Value = 1176
$ws.Range("I22").Value = 450
$ws.Range("K22").Value = 450
$ws.Range("M22").Value = -100
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4624.1816
$ws.Range("I3").Value = 4624.1816
$ws.Range("K3").Value = 13872.5448
$ws.Range("M3").Value = -13760.5448
$ws.Range("H33").Value = 271.44446
$ws.Range("J33").Value = 483
$ws.Range("L33").Value = 2898
$ws.Range("N33").Value = -3464
$ws.Range("H34").Value = 4571.4287
$ws.Range("I34").Value = 6666.6665
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 19999.9995
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -19915.9995
$ws.Range("N34").Value = -9168
$ws.Range("H36").Value = 6210.2
$ws.Range("I36").Value = 262.75
$ws.Range("K36").Value = 788.25
$ws.Range("M36").Value = -619.25
$ws.Range("H68").Value = 4110.6
$ws.Range("J68").Value = 4118.552
$ws.Range("L68").Value = 12355.656
$ws.Range("N68").Value = -13977.656
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""
$ws.Range("H71").Value = 4110.6
$ws.Range("J71").Value = 4118.552
$ws.Range("L71").Value = 37066.96799999999
$ws.Range("N71").Value = -45178.96799999999
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""
$ws.Range("H93").Value = 13800
$ws.Range("J93").Value = 16360
$ws.Range("L93").Value = 49080
$ws.Range("N93").Value = -52824
$ws.Range("H133").Value = 12761.667
$ws.Range("J133").Value = 25000
$ws.Range("L133").Value = 75000
$ws.Range("N133").Value = -85120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5900
$ws.Range("J80").Value = 5900
$ws.Range("L80").Value = 5900
$ws.Range("N80").Value = -7896
$ws.Range("H83").Value = 5900
$ws.Range("J83").Value = 5900
$ws.Range("L83").Value = 29500
$ws.Range("N83").Value = -39484
$ws.Range("H102").Value = 2655.5
$ws.Range("I102").Value = 2348.5
$ws.Range("J102").Value = 2962.5
$ws.Range("K102").Value = 2348.5
$ws.Range("L102").Value = 2962.5
$ws.Range("M102").Value = -726.5
$ws.Range("N102").Value = -6206.5
$ws.Range("H122").Value = 2487.818
$ws.Range("I122").Value = 2777.8333
$ws.Range("K122").Value = 8333.499899999999
$ws.Range("M122").Value = -5883.499899999999
$ws.Range("H123").Value = 50326
$ws.Range("J123").Value = 50326
$ws.Range("L123").Value = 50326
$ws.Range("N123").Value = -55226
$ws.Range("H132").Value = 6453.5
$ws.Range("I132").Value = 6994.25
$ws.Range("K132").Value = 20982.75
$ws.Range("M132").Value = -18452.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6934.273
$ws.Range("I7").Value = 5713.3335
$ws.Range("J7").Value = 8399.4
$ws.Range("K7").Value = 5713.3335
$ws.Range("L7").Value = 8399.4
$ws.Range("M7").Value = -5601.3335
$ws.Range("N7").Value = -8623.4
$ws.Range("H40").Value = 3096.5386
$ws.Range("I40").Value = 2683.2856
$ws.Range("J40").Value = 3578.6667
$ws.Range("K40").Value = 2683.2856
$ws.Range("L40").Value = 3578.6667
$ws.Range("M40").Value = -2547.2856
$ws.Range("N40").Value = -3850.6667
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256
$ws.Range("H93").Value = 1296.5883
$ws.Range("I93").Value = 1296.5883
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1296.5883
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -48.58829999999989
$ws.Range("N93").Value = ""
$ws.Range("H126").Value = 6934.273
$ws.Range("I126").Value = 5713.3335
$ws.Range("J126").Value = 8399.4
$ws.Range("K126").Value = 17140.0005
$ws.Range("L126").Value = 25198.2
$ws.Range("M126").Value = -14670.0005
$ws.Range("N126").Value = -30138.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2216.3635
$ws.Range("I126").Value = 1787.5
$ws.Range("K126").Value = 5362.5
$ws.Range("M126").Value = -2892.5
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = ""
$ws.Range("H132").Value = 1629
$ws.Range("I132").Value = 1898.8572
$ws.Range("K132").Value = 5696.571599999999
$ws.Range("M132").Value = -3166.571599999999
